# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.607.96"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "'1.619.70"
$ws.Range("E3").Value = "  -2.38%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").Value = "'1.005"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "'306.70"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("D7").Value = "'0.3815"
$ws.Range("E7").Value = "  -2.23%  "

$ws.Range("D8").Value = "'0.3745"
$ws.Range("E8").Value = "  -3.16%  "

$ws.Range("D9").Value = "'49.33"
$ws.Range("E9").Value = "  -3.89%  "

$ws.Range("D10").Value = "'1.304"
$ws.Range("E10").Value = "  -4.63%  "

$ws.Range("D11").Value = "'1.006"

$ws.Range("D12").Value = "'0.08227"
$ws.Range("E12").Value = "  -3.29%  "

$ws.Range("D13").Value = "'23.34"
$ws.Range("E13").Value = "  -2.53%  "

$ws.Range("D14").Value = "'6.769"
$ws.Range("E14").Value = "  -6.23%  "

$ws.Range("D15").Value = "'7.615"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("D16").Value = "'0.00001277"
$ws.Range("E16").Value = "  -2.76%  "

$ws.Range("D17").Value = "'1.628.66"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "'92.59"
$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").Value = "'0.06900"
$ws.Range("E19").Value = "  -1.35%  "

$ws.Range("D20").Value = "'18.97"
$ws.Range("E20").Value = "  -4.98%  "

$ws.Range("D21").Value = "'6.764"
$ws.Range("E21").Value = "  -3.04%  "

$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").Value = "'13.35"
$ws.Range("E23").Value = "  -2.31%  "

$ws.Range("D24").Value = "'23.619.02"
$ws.Range("E24").Value = "  -1.66%  "

$ws.Range("D25").Value = "'2.406"
$ws.Range("E25").Value = "  -3.21%  "

$ws.Range("D26").Value = "'2.828"
$ws.Range("E26").Value = "  -8.64%  "

$ws.Range("D27").Value = "'21.52"
$ws.Range("E27").Value = "  -3.26%  "

$ws.Range("D28").Value = "'151.44"
$ws.Range("E28").Value = "  -1.66%  "

$ws.Range("D29").Value = "'5.425"

$ws.Range("D30").Value = "'7.845"
$ws.Range("E30").Value = "  -2.27%  "

$ws.Range("D31").Value = "'133.43"
$ws.Range("E31").Value = "  -4.68%  "

$ws.Range("D32").Value = "'2.482"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("D33").Value = "'1.821.31"
$ws.Range("E33").Value = "  -1.06%  "

$ws.Range("D34").Value = "'0.9657"
$ws.Range("E34").Value = "  -7.83%  "

$ws.Range("D35").Value = "'0.07683"
$ws.Range("E35").Value = "  -5.70%  "

$ws.Range("D36").Value = "'0.02834"
$ws.Range("E36").Value = "  -5.84%  "

$ws.Range("D37").Value = "'6.456"
$ws.Range("E37").Value = "  -4.51%  "

$ws.Range("D38").Value = "'0.2594"
$ws.Range("E38").Value = "  -4.23%  "

$ws.Range("E39").Value = "  -8.00%  "

$ws.Range("D40").Value = "'0.08966"
$ws.Range("E40").Value = "  -2.13%  "

$ws.Range("D41").Value = "'0.7325"
$ws.Range("E41").Value = "  -3.28%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.397"
$ws.Range("E42").Value = "  -1.96%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'13.06"
$ws.Range("E43").Value = "  -4.84%  "

$ws.Range("D44").Value = "'16.28"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").Value = "'0.6751"
$ws.Range("E45").Value = "  -3.95%  "

$ws.Range("D46").Value = "'2.368"
$ws.Range("E46").Value = "  -5.32%  "

$ws.Range("D47").Value = "'4.043"
$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D49").Value = "'0.08096"
$ws.Range("E49").Value = "  -2.47%  "

$ws.Range("D50").Value = "'131.99"
$ws.Range("E50").Value = "  -2.76%  "

$ws.Range("D51").Value = "'1.193"
$ws.Range("E51").Value = "  -3.82%  "
